# bugfix for config file
# - fix misspelled sheet name "materiarDepot" -> "materialDepot"
# - make the (now correctly named) materialDepot sheet the active/selected tab
# - move its selection to F28
# - the previously active sheet ("quarrier") goes back to its default (unselected) state

$wb = $excel.ActiveWorkbook

# Rename the misspelled sheet.
$ws = $wb.Worksheets.Item("materiarDepot")
$ws.Name = "materialDepot"

# Make it the active sheet (flips tabSelected in the sheetViews / activeTab in the workbook).
$ws.Activate()

# Move the selection on the newly active sheet to F28.
$ws.Range("F28").Select()
